$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 195.435389
$ws.Range("H2").Value = 586.306167
$ws.Range("I2").Value = 0.3095741734129938
$ws.Range("J2").Value = 0.3095741734129938
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.1421153333333333
$ws.Range("N2").Value = 0.426346
$ws.Range("Q2").Value = 27.77436545286467
$ws.Range("R2").Value = 249.969289075782
$ws.Range("S2").Value = 0.3095741734129938
$ws.Range("T2").Value = 0.3095741734129938

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 79.82725266666667
$ws.Range("H3").Value = 239.481758
$ws.Range("I3").Value = 0.1264482133280045
$ws.Range("J3").Value = 0.1264482133280045
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.1421153333333333
$ws.Range("N3").Value = 0.426346
$ws.Range("Q3").Value = 11.34467662180756
$ws.Range("R3").Value = 102.102089596268
$ws.Range("S3").Value = 0.1264482133280045
$ws.Range("T3").Value = 0.1264482133280045

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 162.5116576666667
$ws.Range("H4").Value = 487.534973
$ws.Range("I4").Value = 0.2574222219914007
$ws.Range("J4").Value = 0.2574222219914007
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1421153333333333
$ws.Range("N4").Value = 0.426346
$ws.Range("Q4").Value = 23.09539839985089
$ws.Range("R4").Value = 207.858585598658
$ws.Range("S4").Value = 0.2574222219914007
$ws.Range("T4").Value = 0.2574222219914007

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 159.7910413333334
$ws.Range("H5").Value = 479.3731240000001
$ws.Range("I5").Value = 0.2531127028358626
$ws.Range("J5").Value = 0.2531127028358626
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.1421153333333333
$ws.Range("N5").Value = 0.426346
$ws.Range("Q5").Value = 22.70875710276712
$ws.Range("R5").Value = 204.378813924904
$ws.Range("S5").Value = 0.2531127028358626
$ws.Range("T5").Value = 0.2531127028358626

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 33.73857866666666
$ws.Range("H6").Value = 101.215736
$ws.Range("I6").Value = 0.05344268843173843
$ws.Range("J6").Value = 0.05344268843173842
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.1421153333333333
$ws.Range("N6").Value = 0.426346
$ws.Range("Q6").Value = 4.794769353406222
$ws.Range("R6").Value = 43.15292418065599
$ws.Range("S6").Value = 0.05344268843173843
$ws.Range("T6").Value = 0.05344268843173842
